# Append a new application-history entry (row 55) for RBC - Derivative
# Valuation Analyst, mirroring the structure of the existing rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 55

$resumeText = @'

**EDUCATION**
Washington University in St. Louis (WashU)	St. Louis, MO, USA
M.S. in Finance, Quantitative	Sep 2023 - Dec 2024
- GPA: 3.95/4.00, Rank 2/89
- Honors: All-semester Dean's List, Beta Gamma Sigma Award, Charles F. Knight Scholar (Expected)
- Coursework: Continuous-Time Finance (Ph.D. Level), Stochastic Calculus, Exotic & Fixed Income Derivative Pricing, SQL and Database Design, Python & R Machine Learning, LASSO, Quantitative Risk Management, Advanced Corporate Finance I & II
The University of Hong Kong (HKU)	Hong Kong
B.S. in Economics and Finance	Sep 2017 - May 2021
- Rank top 35%, Graduated with 2:1 Distinction, Selected as C.V. Starr Scholarship recipient to exchange at Sciences Po Paris
**PROFESSIONAL EXPERIENCE**
Olin Business School at Washington University in St. Louis	St. Louis, MO, USA
Research Intern	May 2024 - Present
- Selected Project: Python-driven, Real-time Political Speech Sentiment Signal Tool on Abnormal Market Movements
    o	Developed a Selenium Chrome-driver scraper to extract U.S. congressional hearing scripts and related social media accounts
    o	Fine-tuned BERT Topic Model and its sub-models to implement topic classification of the congressional hearing database
    o	Built an anomaly detection model on sentiment trends to flag moments where sentiment sharply deviates from baseline
    o	Created an automated alert mechanism triggered by specific keywords, phrases, or sentiment shifts against fiscal and monetary policies, to enable rapid responses to market-moving and short-term trading opportunities
- Overlaid UK budget announcement shock with market data to visualize how certain political event statements or sentiment trends coincide with market price changes, supporting research project on fiscal shock omitted variable bias (OVB)
Privium Fund Management	Hong Kong
Portfolio Manager - US$ 200M Option Selling (Short Straddle) Strategy	Apr 2022 - Dec 2022
- Co-managed the strategy, applied Black-Scholes and Heston-Nandi GARCH option pricing to optimize premium income
- Ensured macro-overlay and market breadth, evaluated liquidity and volatility metrics to align trades with market trends
- Applied pricing models (Heston, Barra, Black-Litterman, Greek hedging) for underlying asset price prediction and risk alignment
- Controlled portfolio risk through diversification to mitigate contagion and application of risk models to evaluate potential losses
- Optimized strategies, utilized algorithmic trading to minimize slippage, enhance efficiency, and reduce risk exposure
- Integrated transaction costs, stress testing, and market impact analysis to guarantee accurate performance metrics
Yong Rong Asset	Hong Kong
Junior Trader - US$ 30M High Conviction Sub-fund	Jun 2021 - Jan 2022
- Analyzed and constructed trade flows, submitting orders and collaborating with counterparties to optimize trading outcomes
- Utilized strategies such as limit orders, VWAP, TWAP to achieve optimal execution and transaction cost control (TCA)
- Supported strategic asset allocation, used VBA to streamline trading script maintenance and Net Asset Value report process
Peak Global Investments	Hong Kong
Private Equity Intern	Sep 2020 - May 2021
- Researched cryptocurrency exchanges across Asia and Europe, liaised with senior executives to prepare for acquisitions
- Assessed client pain points and formulated pitch decks for management team in business development presentations
- Collaborated with world's largest crypto exchange, applied API to assess targets' trading volumes and their authenticity
**PROJECT EXPERIENCE**
Microstructure-Informed End-of-Day Frequency Trading Strategy	St. Louis, MO, USA
Trading Strategy Design	Oct 2024 - Present
- Currently developing a proprietary end-of-day (EoD) frequency trading strategy, leveraging market microstructure analysis and machine learning techniques (LASSO, LSTM, kNN) to predict price movements and execute trades in the U.S. equities market
Hull-White Model Calibration for At-the-Money (ATM) Caplets and Caps	St. Louis, MO, USA
Outputs: fredhli.github.io/projects/#hull-white-model-calibration  	Mar 2024 - May 2024
- Implemented closed-form caplet pricing solutions and Monte-Carlo price check, to optimize Hull-White model parameters against ATM Caplet market data, achieving high accuracy in long-maturity cap pricing > 15 years with less than 3% function value loss
**SKILLS AND QUALIFICATIONS**
Certificates 	  CFA Level I, HKSFC Type 4 & 9: Advise on Securities & Asset Management License
Teaching	  TA for: Options, Futures and Derivative Securities (Undergraduate); Behavioral Finance (Graduate)
Research	  RA for: Research on PEVC-backed companies - under Prof. Minmo Gahng, Cornell University
Volunteering	  NGO Marketing Director, Soap Cycling HKU; Village School Teacher, Beyond the Pivot HKU
Programming	  Proficient in Python, R, SQL, Git, VBA, LaTeX; Intermediate in Stata, MATLAB; Basic in JavaScript
Work Permits	  Hong Kong SAR (Permanent Citizenship), Canada (OWP with Citizenship Assurance), USA (OPT)

'@

$jobDescText = @'

About the job
Job Summary
Job Description
What is the opportunity?
The Derivative Valuation Analyst is responsible for the timely and accurate analysis, valuation, validation and reconciliation of Listed, OTC (Over the Counter) Bilateral and OTC Cleared Derivatives, as well as of alternative investment products, impacting key lines of business including fund services, pension reporting, insurance reporting, custody reporting and all value added services.
Please note, this position may cover a late and/or rotating shift schedule.
What will you do?
Perform the daily
Routine tasks such as setup of Listed Derivatives and OTC Derivatives, in Valuation systems
Tasks for listed Derivatives pricing and analyzing price movements
OTC/alternative instruments valuation and measure impact of the daily changes in market data inputs to value of the instrument
Reconciliation against third party statements
Act as the main liaison with Investment Managers and the internal Trades processing team on new deals and supports the communication of deal price variances that exceed tolerances to the Investment Managers and obtain confirmation
Participate in research of appropriate valuation methodology for new complex OTC products traded by clients
Participate in the execution around key project deliverables and identify opportunities to stream line processes and increase productivity
Ensure compliance with all policies, procedures and standards for all aspects of the business as defined through Risk, Compliance, Operational reviews and audits
Cover early/late shifts to support the business within global operating model
What do you need to succeed?
Must-have:
Undergraduate Degree in Finance, Math or Engineering
2+ years of solid experience with Derivatives in a middle or back office environment
Experience with OTC Derivatives valuation systems and Workflows
Strong understanding of Reuters, Bloomberg and other financial data providers
Strong PC skills (advanced proficiency in MS Excel, some knowledge of basic coding such as VBA)
What's in it for you?
We thrive on the challenge to be our best, progressive thinking to keep growing, and working together to deliver trusted advice to help our clients thrive and communities prosper. We care about each other, reaching our potential, making a difference to our communities, and achieving success that is mutual.
Excellent exposure to communicate with various business partners and stakeholders in Investor Services and within other platforms as appropriate
Opportunity to obtain hands-on experience throughout your role
Working with an exciting, close-knit, supportive & dynamic group
Opportunity to collaborate with other business segments within the bank
Excellent career development and progression opportunities
A comprehensive Total Rewards Program including bonuses and flexible benefits
Competitive compensation
Job Skills
Active Learning, Adaptability, Business Appraisals, Critical Thinking, Customer Service, Decision Making, Effectiveness Measurement, Operational Delivery, Process Improvements
Additional Job Details
Address:
RBC CENTRE, 155 WELLINGTON ST W:TORONTO
City:
TORONTO
Country:
Canada
Work hours/week:
37.5
Employment Type:
Full time
Platform:
WEALTH MANAGEMENT
Job Type:
Regular
Pay Type:
Salaried
Posted Date:
2024-10-17
Application Deadline:
2024-11-27
Note: Applications will be accepted until 11:59 PM on the day prior to the application deadline date above
Inclusion and Equal Opportunity Employment
At RBC, we embrace diversity and inclusion for innovation and growth. We are committed to building inclusive teams and an equitable workplace for our employees to bring their true selves to work. We are taking actions to tackle issues of inequity and systemic bias to support our diverse talent, clients and communities.
We also strive to provide an accessible candidate experience for our prospective employees with different abilities. Please let us know if you need any accommodations during the recruitment process.
Join our Talent Community
Stay in-the-know about great career opportunities at RBC. Sign up and get customized info on our latest jobs, career tips and Recruitment events that matter to you.
Expand your limits and create a new future together at RBC. Find out how we use our passion and drive to enhance the well-being of our clients and communities

'@

$ws.Cells.Item($row, 1).Value = "November 20, 2024"
$ws.Cells.Item($row, 2).Value = "RBC"
$ws.Cells.Item($row, 3).Value = "Derivative Valuation Analyst"
$ws.Cells.Item($row, 4).Value = "Cover Letter - Fred Li.pdf"
$ws.Cells.Item($row, 5).Value = $resumeText
$ws.Cells.Item($row, 6).Value = "Toronto"
$ws.Cells.Item($row, 7).Value = $jobDescText
